$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.03644326811879168
$ws.Range("E2").Value = 0.393865553482251

$ws.Range("D3").Value = 0.2808990339019236
$ws.Range("E3").Value = 6.57676811390915

$ws.Range("D4").Value = -0.0001415336340117438
$ws.Range("E4").Value = 0.01083971667409135

$ws.Range("D5").Value = -0.3671307302191725
$ws.Range("E5").Value = 6.128582539354844

$ws.Range("D6").Value = 0.09342891388835872
$ws.Range("E6").Value = 0.4583074153935279

$ws.Range("D7").Value = 0.5989335063986924
$ws.Range("E7").Value = 6.630267284365771

$ws.Range("D8").Value = -0.0001555332065405596
$ws.Range("E8").Value = 0.007012150680553096

$ws.Range("D9").Value = -0.4950569817774801
$ws.Range("E9").Value = 6.199878703985113
